$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for "Puerro" needs to be inserted right after the
# existing row 3, which pushes every subsequent record (previously rows 4-95)
# down by one row (to rows 5-96), growing the used range to A1:R96.
$ws.Rows.Item(4).Insert()

# The descriptive columns (market/region/category/etc.) are identical for
# every record in this sheet, so copy them from the row directly above.
$ws.Range("A4").Value = $ws.Range("A3").Value()
$ws.Range("B4").Value = $ws.Range("B3").Value()
$ws.Range("C4").Value = $ws.Range("C3").Value()
$ws.Range("E4").Value = $ws.Range("E3").Value()
$ws.Range("F4").Value = $ws.Range("F3").Value()
$ws.Range("G4").Value = $ws.Range("G3").Value()
$ws.Range("H4").Value = $ws.Range("H3").Value()
$ws.Range("I4").Value = $ws.Range("I3").Value()
$ws.Range("N4").Value = $ws.Range("N3").Value()
$ws.Range("O4").Value = $ws.Range("O3").Value()
$ws.Range("Q4").Value = $ws.Range("Q3").Value()
$ws.Range("R4").Value = $ws.Range("R3").Value()

# New record's own data: date, volume, min/max/avg price, $/Kg.
$ws.Range("D4").Value = 44699
$ws.Range("J4").Value = 160
$ws.Range("K4").Value = 7000
$ws.Range("L4").Value = 8000
$ws.Range("M4").Value = 7500
$ws.Range("P4").Value = 375
